# print-template.docx edit:
#  1. Heading placeholder "{Month} {Year}" -> "{date}" (split across 3 runs:
#     "{" / "d" / "ate}" to mirror the captured edit's run shape).
#  2. Table cell placeholder text stays "{#entries}{Date}" but is re-split
#     into "{#" / "e" / "ntries}{Date}" runs, and the stray gramStart/gramEnd
#     w:proofErr markers around the old run boundary are removed.
#  3. The header's inline picture run gets an explicit <w:noProof/> so Word
#     doesn't try to spell/grammar-check the drawing run.

$d = $word.ActiveDocument

# --- 1: heading "{Month} {Year}" -> "{date}" --------------------------------
$headingPara = $d.Paragraphs(1).Range
$headingPara.End = $headingPara.End - 1          # exclude the paragraph mark
if ($headingPara.Text -eq "{Month} {Year}") {
    $headingPara.Text = "{date}"
}

$headingStart = $d.Paragraphs(1).Range.Start
# Force Word to break "{date}" into separate runs at the "d" boundary by
# toggling a character property on just that letter (bold on, then back off).
$dChar = $d.Range($headingStart + 1, $headingStart + 2)
if ($dChar.Text -eq "d") {
    $dChar.Bold = 1
    $dChar.Bold = 0
}

# --- 2: "{#entries}{Date}" cell - re-split runs & drop proofErr marks -------
$entriesPara = $d.Paragraphs(7).Range
$entriesBase = $entriesPara.Start
$entriesCheck = $d.Range($entriesBase, $entriesBase + 16)
if ($entriesCheck.Text -eq "{#entries}{Date}") {
    # Rewrite through a scratch value first so the identical-text case still
    # collapses to one freshly minted run (with the proofErr marks dropped),
    # then restore the real text into that single run.
    $scratch = $d.Range($entriesBase, $entriesBase + 16)
    $scratch.Text = "................"
    $restored = $d.Range($entriesBase, $entriesBase + 16)
    $restored.Text = "{#entries}{Date}"
}

# Re-split "{#entries}{Date}" after the first "e" (i.e. "{#" / "e" /
# "ntries}{Date}") the same way as above: toggle bold on just that letter.
$eChar = $d.Range($entriesBase + 2, $entriesBase + 3)
if ($eChar.Text -eq "e") {
    $eChar.Bold = 1
    $eChar.Bold = 0
}

# --- 3: header picture run gets <w:noProof/> --------------------------------
$section = $d.Sections(1)
$headers = $section.Headers
for ($i = 1; $i -le $headers.Count; $i++) {
    $hf = $headers.Item($i)
    if ($hf.Exists -and $hf.Range.InlineShapes.Count -gt 0) {
        $pic = $hf.Range.InlineShapes.Item(1)
        $pic.Range.NoProofing = 1
    }
}

Write-Output "Heading: [$($d.Paragraphs(1).Range.Text)]"
Write-Output "Entries: [$($d.Paragraphs(7).Range.Text)]"
